$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

$ws.Range("B47").Value = "SingleUseId55"
$ws.Range("C47").Value = "Default"
$ws.Range("D47").Value = "Left"
$ws.Range("E47").Value = "LTR"
$ws.Range("F47").Value = "THRESHOLD"

$ws.Range("B48").Value = "SingleUseId56"
$ws.Range("C48").Value = "Default"
$ws.Range("D48").Value = "Left"
$ws.Range("E48").Value = "LTR"
$ws.Range("F48").Value = "SLOPE"

$ws.Range("B49").Value = "SingleUseId57"
$ws.Range("C49").Value = "Default"
$ws.Range("D49").Value = "Center"
$ws.Range("E49").Value = "LTR"
$ws.Range("F49").Value = "Detect"

$ws.Range("B50").Value = "SingleUseId58"
$ws.Range("C50").Value = "Default"
$ws.Range("D50").Value = "Center"
$ws.Range("E50").Value = "LTR"
$ws.Range("F50").Value = "<value> mV"

$ws.Range("B51").Value = "SingleUseId59"
$ws.Range("C51").Value = "Default"
$ws.Range("D51").Value = "Left"
$ws.Range("E51").Value = "LTR"
$ws.Range("F51").NumberFormat = "@"
$ws.Range("F51").Value = "0"
$ws.Range("F51").Style = "Normal"

$ws.Range("B52").Value = "SingleUseId62"
$ws.Range("C52").Value = "Default"
$ws.Range("D52").Value = "Left"
$ws.Range("E52").Value = "LTR"
$ws.Range("F52").Value = "External"

$ws.Range("B53").Value = "SingleUseId63"
$ws.Range("C53").Value = "Default"
$ws.Range("D53").Value = "Center"
$ws.Range("E53").Value = "LTR"
$ws.Range("F53").Value = "Internal`nRubid"

$ws.Range("B54").Value = "SingleUseId64"
$ws.Range("C54").Value = "Default"
$ws.Range("D54").Value = "Center"
$ws.Range("E54").Value = "LTR"
$ws.Range("F54").Value = "Internal`nQuartz"

$ws.Range("B55").Value = "SingleUseId65"
$ws.Range("C55").Value = "Default"
$ws.Range("D55").Value = "Center"
$ws.Range("E55").Value = "LTR"
$ws.Range("F55").Value = "<value>"

$ws.Range("B56").Value = "SingleUseId66"
$ws.Range("C56").Value = "Default"
$ws.Range("D56").Value = "Left"
$ws.Range("E56").Value = "LTR"
$ws.Range("F56").Value = "HF INPUT"

$ws.Range("B57").Value = "SingleUseId67"
$ws.Range("C57").Value = "Large"
$ws.Range("D57").Value = "Left"
$ws.Range("E57").Value = "LTR"
$ws.Range("F57").Value = "ON"

$ws.Range("B58").Value = "SingleUseId68"
$ws.Range("C58").Value = "Large"
$ws.Range("D58").Value = "Left"
$ws.Range("E58").Value = "LTR"
$ws.Range("F58").Value = "OFF"

$ws.Range("B59").Value = "SingleUseId69"
$ws.Range("C59").Value = "Default"
$ws.Range("D59").Value = "Left"
$ws.Range("E59").Value = "LTR"
$ws.Range("F59").Value = "GATE"

$ws.Range("B60").Value = "SingleUseId70"
$ws.Range("C60").Value = "Default"
$ws.Range("D60").Value = "Center"
$ws.Range("E60").Value = "LTR"
$ws.Range("F60").Value = "<value> ms"

$ws.Range("B61").Value = "SingleUseId71"
$ws.Range("C61").Value = "Default"
$ws.Range("D61").Value = "Left"
$ws.Range("E61").Value = "LTR"
$ws.Range("F61").Value = "MES SETUP"

$ws.Range("B62").Value = "SingleUseId73"
$ws.Range("C62").Value = "Default"
$ws.Range("D62").Value = "Center"
$ws.Range("E62").Value = "LTR"
$ws.Range("F62").Value = "<value>"

$ws.Range("B63").Value = "SingleUseId74"
$ws.Range("C63").Value = "Default"
$ws.Range("D63").Value = "Left"
$ws.Range("E63").Value = "LTR"
$ws.Range("F63").NumberFormat = "@"
$ws.Range("F63").Value = "0"
$ws.Range("F63").Style = "Normal"

$ws.Range("B64").Value = "SingleUseId75"
$ws.Range("C64").Value = "Default"
$ws.Range("D64").Value = "Left"
$ws.Range("E64").Value = "LTR"
$ws.Range("F64").Value = "Stamps Number"

$ws.Range("B65").Value = "SingleUseId76"
$ws.Range("C65").Value = "Large"
$ws.Range("D65").Value = "Left"
$ws.Range("E65").Value = "LTR"
$ws.Range("F65").Value = "X"

$ws.Range("B66").Value = "SingleUseId77"
$ws.Range("C66").Value = "Default"
$ws.Range("D66").Value = "Left"
$ws.Range("E66").Value = "LTR"
$ws.Range("F66").Value = "Repeat"

$ws.Range("B67").Value = "SingleUseId78"
$ws.Range("C67").Value = "Default"
$ws.Range("D67").Value = "Center"
$ws.Range("E67").Value = "LTR"
$ws.Range("F67").Value = "<value>"

$ws.Range("B68").Value = "SingleUseId79"
$ws.Range("C68").Value = "Default"
$ws.Range("D68").Value = "Left"
$ws.Range("E68").Value = "LTR"
$ws.Range("F68").NumberFormat = "@"
$ws.Range("F68").Value = "0"
$ws.Range("F68").Style = "Normal"

$ws.Range("B69").Value = "SingleUseId84"
$ws.Range("C69").Value = "Default"
$ws.Range("D69").Value = "Left"
$ws.Range("E69").Value = "LTR"
$ws.Range("F69").Value = "Single"

$ws.Range("B70").Value = "SingleUseId85"
$ws.Range("C70").Value = "Default"
$ws.Range("D70").Value = "Center"
$ws.Range("E70").Value = "LTR"
$ws.Range("F70").Value = "Continuous"

$ws.Range("B71").Value = "SingleUseId86"
$ws.Range("C71").Value = "Default"
$ws.Range("D71").Value = "Center"
$ws.Range("E71").Value = "LTR"
$ws.Range("F71").Value = "Stamps"

$ws.Range("B72").Value = "SingleUseId87"
$ws.Range("C72").Value = "Default"
$ws.Range("D72").Value = "Center"
$ws.Range("E72").Value = "LTR"
$ws.Range("F72").Value = "Single"

$ws.Range("B73").Value = "SingleUseId88"
$ws.Range("C73").Value = "Default"
$ws.Range("D73").Value = "Center"
$ws.Range("E73").Value = "LTR"
$ws.Range("F73").Value = "Continuous"

$ws.Range("B74").Value = "SingleUseId89"
$ws.Range("C74").Value = "Default"
$ws.Range("D74").Value = "Center"
$ws.Range("E74").Value = "LTR"
$ws.Range("F74").Value = "Stamps"
